$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18:F18").NumberFormat = "@"

$ws.Range("A18").Value = "Totals"
$ws.Range("B18").Value = ""
$ws.Range("C18").Value = ""
$ws.Range("D18").Value = "13790687.00"
$ws.Range("E18").Value = "0"
$ws.Range("F18").Value = "13790687.00"
